# PrezConfig1.xlsx edit: add an "Excluido" column to the "Recursos" sheet
# (mirroring the already-existing "Excluido" column on "Usuarios"), make
# "Usuarios" the active / selected sheet, and refresh the column widths /
# page margins that Excel rewrites whenever the workbook is resaved.

$wb = $excel.ActiveWorkbook
$wsUsuarios = $wb.Worksheets.Item("Usuarios")
$wsRecursos = $wb.Worksheets.Item("Recursos")

# --- Recursos: add the "Excluido" header in column C -----------------
$wsRecursos.Range("C1").Value = "Excluido"

# Copy the header formatting (bold, centered, bordered) from B1 onto the
# new C1 header cell so it matches A1/B1.
$wsRecursos.Range("B1").Copy()
$wsRecursos.Range("C1").PasteSpecial(-4122)

# Column widths: Usuarios!A fits the longest user name, Recursos!A:B fit
# the longest course/context string.
$wsUsuarios.Columns.Item(1).ColumnWidth = 30.1
$wsRecursos.Columns.Item(1).ColumnWidth = 47.833333333333333
$wsRecursos.Columns.Item(2).ColumnWidth = 47.833333333333333

# Page margins: reset to Excel's normal defaults (0.7"/0.7"/0.75"/0.75"/0.3"/0.3").
foreach ($ws in @($wsUsuarios, $wsRecursos)) {
    $ws.PageSetup.LeftMargin = 0.7 * 72
    $ws.PageSetup.RightMargin = 0.7 * 72
    $ws.PageSetup.TopMargin = 0.75 * 72
    $ws.PageSetup.BottomMargin = 0.75 * 72
    $ws.PageSetup.HeaderMargin = 0.3 * 72
    $ws.PageSetup.FooterMargin = 0.3 * 72
}

# --- Active sheet / selection: Usuarios becomes the visible tab ------
$wsUsuarios.Activate()
$wsUsuarios.Range("A2:A14").Select()
